$wb = $excel.ActiveWorkbook

# This script applies per-cell numeric updates to the Leve profit-tracking tables
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, refreshing computed market-price
# and profit columns (H:N) as produced by the scheduled pricing-data runner.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1869.2128
$ws.Range("I40").Value = 1871.3684
$ws.Range("J40").Value = 1860.1111
$ws.Range("K40").Value = 1871.3684
$ws.Range("L40").Value = 1860.1111
$ws.Range("M40").Value = -1696.3684
$ws.Range("N40").Value = -2210.1111

$ws.Range("H45").Value = 6995
$ws.Range("I45").Value = 7000
$ws.Range("J45").Value = 6990
$ws.Range("K45").Value = 21000
$ws.Range("L45").Value = 20970
$ws.Range("M45").Value = -20808
$ws.Range("N45").Value = -21354

$ws.Range("H99").Value = 573.1429000000001
$ws.Range("J99").Value = 1000
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996

$ws.Range("H132").Value = 2331.5652
$ws.Range("I132").Value = 2423.6365
$ws.Range("J132").Value = 306
$ws.Range("K132").Value = 7270.9095
$ws.Range("L132").Value = 918
$ws.Range("M132").Value = -4740.9095
$ws.Range("N132").Value = -5978

$ws.Range("H133").Value = 20998.572
$ws.Range("J133").Value = 20998.572
$ws.Range("L133").Value = 20998.572
$ws.Range("N133").Value = -31118.572

$ws.Range("H138").Value = 3457.2698
$ws.Range("I138").Value = 1601.1364
$ws.Range("J138").Value = 4453.244
$ws.Range("K138").Value = 4803.4092
$ws.Range("L138").Value = 13359.732
$ws.Range("M138").Value = 336.5907999999999
$ws.Range("N138").Value = -23639.732


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 7283.3335
$ws.Range("I28").Value = 7283.3335
$ws.Range("K28").Value = 7283.3335
$ws.Range("M28").Value = -7091.3335

$ws.Range("H32").Value = 7707.906
$ws.Range("I32").Value = 8830.852999999999
$ws.Range("J32").Value = 5698.421
$ws.Range("K32").Value = 8830.852999999999
$ws.Range("L32").Value = 5698.421
$ws.Range("M32").Value = -8543.852999999999
$ws.Range("N32").Value = -6272.421

$ws.Range("H99").Value = 7283.3335
$ws.Range("I99").Value = 7283.3335
$ws.Range("K99").Value = 7283.3335
$ws.Range("M99").Value = -4288.3335

$ws.Range("H109").Value = 38500
$ws.Range("J109").Value = 38500
$ws.Range("L109").Value = 38500
$ws.Range("N109").Value = -41274

$ws.Range("H114").Value = 129999.75
$ws.Range("J114").Value = 129999.75
$ws.Range("L114").Value = 129999.75
$ws.Range("N114").Value = -138677.75

$ws.Range("H119").Value = 29698
$ws.Range("J119").Value = 29698
$ws.Range("L119").Value = 29698
$ws.Range("N119").Value = -39374

$ws.Range("H132").Value = 2213.5425
$ws.Range("I132").Value = 1276.0769
$ws.Range("J132").Value = 4041.6
$ws.Range("K132").Value = 3828.2307
$ws.Range("L132").Value = 12124.8
$ws.Range("M132").Value = -1298.2307
$ws.Range("N132").Value = -17184.8


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H118").Value = 10307346
$ws.Range("J118").Value = 10307346
$ws.Range("L118").Value = 10307346
$ws.Range("N118").Value = -10310660


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 301.33334
$ws.Range("I22").Value = 224
$ws.Range("K22").Value = 224
$ws.Range("M22").Value = 126

$ws.Range("H62").Value = 3809
$ws.Range("I62").Value = 3638.8572
$ws.Range("K62").Value = 3638.8572
$ws.Range("M62").Value = -3014.8572

$ws.Range("H65").Value = 3809
$ws.Range("I65").Value = 3638.8572
$ws.Range("K65").Value = 18194.286
$ws.Range("M65").Value = -15074.286

$ws.Range("H132").Value = 3388.5334
$ws.Range("I132").Value = 2546.2222
$ws.Range("K132").Value = 7638.6666
$ws.Range("M132").Value = -5108.6666


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1488.4286
$ws.Range("I2").Value = 2532.25
$ws.Range("J2").Value = 96.666664
$ws.Range("K2").Value = 15193.5
$ws.Range("L2").Value = 579.999984
$ws.Range("M2").Value = -15080.5
$ws.Range("N2").Value = -805.999984

$ws.Range("H7").Value = 941.8333
$ws.Range("I7").Value = 333.33334
$ws.Range("J7").Value = 1144.6666
$ws.Range("K7").Value = 1000.00002
$ws.Range("L7").Value = 3433.9998
$ws.Range("M7").Value = -888.0000200000001
$ws.Range("N7").Value = -3657.9998

$ws.Range("H63").Value = 3000
$ws.Range("J63").Value = 5666.6665
$ws.Range("L63").Value = 16999.9995
$ws.Range("N63").Value = -18497.9995

$ws.Range("H66").Value = 3000
$ws.Range("J66").Value = 5666.6665
$ws.Range("L66").Value = 50999.9985
$ws.Range("N66").Value = -58487.9985

$ws.Range("H68").Value = 2598.9875
$ws.Range("I68").Value = 3021.7556
$ws.Range("K68").Value = 9065.266799999999
$ws.Range("M68").Value = -8254.266799999999

$ws.Range("H71").Value = 2598.9875
$ws.Range("I71").Value = 3021.7556
$ws.Range("K71").Value = 27195.8004
$ws.Range("M71").Value = -23139.8004

$ws.Range("H107").Value = 1161.5883
$ws.Range("I107").Value = 343.77777
$ws.Range("K107").Value = 1031.33331
$ws.Range("M107").Value = 888.66669

$ws.Range("H129").Value = 12822584
$ws.Range("I129").Value = 37038252
$ws.Range("J129").Value = 2525.9412
$ws.Range("K129").Value = 111114756
$ws.Range("L129").Value = 7577.823600000001
$ws.Range("M129").Value = -111109756
$ws.Range("N129").Value = -17577.8236

$ws.Range("H131").Value = 11703365
$ws.Range("I131").Value = 6250503.5
$ws.Range("J131").Value = 12821901
$ws.Range("K131").Value = 18751510.5
$ws.Range("L131").Value = 38465703
$ws.Range("M131").Value = -18746470.5
$ws.Range("N131").Value = -38475783

$ws.Range("H138").Value = 9059
$ws.Range("I138").Value = 8620
$ws.Range("J138").Value = 10083.333
$ws.Range("K138").Value = 25860
$ws.Range("L138").Value = 30249.999
$ws.Range("M138").Value = -20720
$ws.Range("N138").Value = -40529.999

$ws.Range("H139").Value = 4768.081
$ws.Range("I139").Value = 7323.0625
$ws.Range("J139").Value = 2821.4285
$ws.Range("K139").Value = 21969.1875
$ws.Range("L139").Value = 8464.2855
$ws.Range("M139").Value = -16829.1875
$ws.Range("N139").Value = -18744.2855

$ws.Range("H140").Value = 6937.625
$ws.Range("I140").Value = 6937.625
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 20812.875
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -15632.875
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 17660.45
$ws.Range("I141").Value = 14356.556
$ws.Range("J141").Value = 20363.637
$ws.Range("K141").Value = 43069.66800000001
$ws.Range("L141").Value = 61090.91099999999
$ws.Range("M141").Value = -37889.66800000001
$ws.Range("N141").Value = -71450.91099999999


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6168.8726
$ws.Range("I136").Value = 4710
$ws.Range("J136").Value = 9167.666999999999
$ws.Range("K136").Value = 14130
$ws.Range("L136").Value = 27503.001
$ws.Range("M136").Value = -11580
$ws.Range("N136").Value = -32603.001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 54050.4
$ws.Range("J27").Value = 54050.4
$ws.Range("L27").Value = 54050.4
$ws.Range("N27").Value = -54188.4

$ws.Range("H126").Value = 1125.8889
$ws.Range("I126").Value = 729.4167
$ws.Range("J126").Value = 1918.8334
$ws.Range("K126").Value = 2188.2501
$ws.Range("L126").Value = 5756.5002
$ws.Range("M126").Value = 281.7498999999998
$ws.Range("N126").Value = -10696.5002

$ws.Range("H130").Value = 34000
$ws.Range("J130").Value = 34000
$ws.Range("L130").Value = 34000
$ws.Range("N130").Value = -44040

$ws.Range("H132").Value = 1942.7234
$ws.Range("I132").Value = 1391.9656
$ws.Range("J132").Value = 2830.0557
$ws.Range("K132").Value = 4175.8968
$ws.Range("L132").Value = 8490.167099999999
$ws.Range("M132").Value = -1645.8968
$ws.Range("N132").Value = -13550.1671

$ws.Range("H136").Value = 1322.28
$ws.Range("I136").Value = 881.6842
$ws.Range("J136").Value = 2717.5
$ws.Range("K136").Value = 2645.0526
$ws.Range("L136").Value = 8152.5
$ws.Range("M136").Value = -95.05259999999998
$ws.Range("N136").Value = -13252.5

